# Added the Quad area function.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - squared differences between point 2 (row4) and point 1 (row2), plus distance
$ws.Range("F6").Formula = "=(B4-B2)^2"
$ws.Range("G6").Formula = "=(C4-C2)^2"
$ws.Range("H6").Formula = "=(D4-D2)^2"
$ws.Range("J6").Formula = "=SQRT(SUM(F6:H6) )"

# Row 28 - diagonal vector 1-3
$ws.Range("B28").Formula = "=B1-B3"
$ws.Range("C28").Formula = "=C1-C3"
$ws.Range("D28").Formula = "=D1-D3"

# Row 29 - diagonal vector 2-4
$ws.Range("B29").Formula = "=B2-B4"
$ws.Range("C29").Formula = "=C2-C4"
$ws.Range("D29").Formula = "=D2-D4"

# Row 32 - cross product of diagonals and quad area
$ws.Range("B32").Formula = "=C28*D29-D28*C29"
$ws.Range("C32").Formula = "=D28*B29-B28*D29"
$ws.Range("D32").Formula = "=B28*C29-C28*B29"
$ws.Range("F32").Formula = "=SQRT(B32^2+C32^2+D32^2)/2"

# Row 34 - edge vector 1-2, dot product with normal (planarity check)
$ws.Range("B34").Formula = "=B1-B2"
$ws.Range("C34").Formula = "=C1-C2"
$ws.Range("D34").Formula = "=D1-D2"
$ws.Range("F34").Formula = "=B34*B32+C34*C32+D34*D32"

# Row 36 - edge vector 3-4, dot product with normal (planarity check)
$ws.Range("B36").Formula = "=B3-B4"
$ws.Range("C36").Formula = "=C3-C4"
$ws.Range("D36").Formula = "=D3-D4"
$ws.Range("F36").Formula = "=B36*B32+C36*C32+D36*D32"

# Update the active selection to match the new working cell
$ws.Range("F32").Select()
